$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first column header from "ID" to "S NO"
$ws.Range("A1").Value = "S NO"
